$d = $word.ActiveDocument

# --- Step 1: the lone _GoBack bookmark in the original document sits at
# the end of the "Alienigena -> Ataques psiquicos..." paragraph (in the
# Characters section, near the end of the body). In the target revision
# it has moved to the end of the new "Interfaz..." paragraph added below.
# Remove it now, while it is still the *only* bookmark of that name, so
# the later by-name lookup/deletion can't be ambiguous once the new one
# (carried inside the inserted fragment) exists too.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: insert the new "things to discuss today" paragraph block
# right after the first body paragraph of the "Introduccion" section
# (i.e. right before the page-break paragraph that starts the "Story"
# section). We append an empty paragraph mark after paragraph 2 and then
# replace that new (still-empty) paragraph's contents with the OOXML
# fragment for the whole block - this reproduces the exact run / spell
# -check proofErr / bookmark structure of the target revision instead of
# just the plain text.
$anchor = $d.Paragraphs(2)
$anchor.Range.InsertParagraphAfter()
$target = $d.Paragraphs(3)

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>Numero de mapas para el inicio</w:t></w:r></w:p><w:p><w:r><w:t>Numero de niveles por mapa</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Tres ramas y tres clases (humano/bestia/hibrido – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ranged</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mele</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-tanque)</w:t></w:r></w:p><w:p><w:r><w:t>Historia personajes</w:t></w:r></w:p><w:p><w:r><w:t>Armas de cada uno</w:t></w:r></w:p><w:p><w:r><w:t>Enemigos</w:t></w:r></w:p><w:p><w:r><w:t>Boses</w:t></w:r></w:p><w:p><w:r><w:t>Tienda</w:t></w:r></w:p><w:p><w:r><w:t>Pickups</w:t></w:r></w:p><w:p><w:r><w:t>Ítems por nivel</w:t></w:r></w:p><w:p><w:r><w:t>Secretos (salas secretas)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Alguna mecánica </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (es decir tener otra opción y no solo disparar, véase bombas y activas en Isaac roll y </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>blanks</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>enter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>the</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gungeon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>Transporte rápido entre salas (puesto que ahora es lineal el mapa puede ser un coñazo volver a coger cosas)</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Walljump</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>?</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Menu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> inicio</w:t></w:r></w:p><w:p><w:r><w:t>Menú pausa</w:t></w:r></w:p><w:p><w:r><w:t>Interfaz (que queremos que se vea y como, ejemplo: vida)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.Range.InsertXML($xml)

# The fragment above ends with an empty paragraph (right after
# "Interfaz..."), but a *trailing* empty paragraph at the very end of an
# InsertXML fragment does not materialize as a new paragraph mark here.
# Locate the "Interfaz..." paragraph we just inserted and append the
# missing trailing empty paragraph explicitly.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Interfaz (que queremos")) {
        $d.Paragraphs($i).Range.InsertParagraphAfter()
        break
    }
}
